# Update the "Score" worksheet:
#  - Populate new columns H (=15), I (=0), J, K for rows 2-21.
#  - Add column L = J + K*20 (row 2 as a standalone formula, rows 3-21 as a
#    shared formula, matching the existing pattern used by column G).
#  - Clear the stray P16 label (this also lets Excel drop the now-unused
#    ";" shared string when it re-saves the workbook).
#  - Update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")
$ws.Activate()

# New data for columns H, I, J, K (row => H, I, J, K)
$data = @(
    @{Row=2;  H=15; I=0; J=10320; K=36}
    @{Row=3;  H=15; I=0; J=9055;  K=30}
    @{Row=4;  H=15; I=0; J=9380;  K=57}
    @{Row=5;  H=15; I=0; J=8280;  K=483}
    @{Row=6;  H=15; I=0; J=7180;  K=15}
    @{Row=7;  H=15; I=0; J=9640;  K=44}
    @{Row=8;  H=15; I=0; J=8220;  K=109}
    @{Row=9;  H=15; I=0; J=8900;  K=181}
    @{Row=10; H=15; I=0; J=11810; K=79}
    @{Row=11; H=15; I=0; J=8830;  K=33}
    @{Row=12; H=15; I=0; J=9980;  K=15}
    @{Row=13; H=15; I=0; J=10300; K=96}
    @{Row=14; H=15; I=0; J=10220; K=32}
    @{Row=15; H=15; I=0; J=6490;  K=29}
    @{Row=16; H=15; I=0; J=10200; K=90}
    @{Row=17; H=15; I=0; J=8900;  K=102}
    @{Row=18; H=15; I=0; J=7910;  K=32}
    @{Row=19; H=15; I=0; J=16020; K=57}
    @{Row=20; H=15; I=0; J=7250;  K=17}
    @{Row=21; H=15; I=0; J=8510;  K=126}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
    $ws.Cells.Item($r, 11).Value = $entry.K
}

# Column L: J + K*20. Row 2 gets its own formula, rows 3-21 share one.
$ws.Range("L2").Formula = "=J2+K2*20"
$ws.Range("L3:L21").Formula = "=J3+K3*20"

# Remove the old, now-unused "P16" label cell (it referenced the stray ";"
# shared string that is no longer needed anywhere in the workbook).
$ws.Cells.Item(16, 16).ClearContents()

# Update the selection to match the saved workbook state.
$ws.Range("J22").Select()
